$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are plain text; set directly.
# Column D (Price) values must be forced to Text format to preserve exact
# string formatting (e.g. trailing zeros, thousand-dot separators) because
# Excel would otherwise auto-convert them to numbers and lose formatting.
# Column E (Volume) values contain leading/trailing spaces and stay text naturally.

$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "53.451.18"
$ws.Range("E2").Value = "  -11.94%  "
$ws.Range("D3").Value = "2.329.18"
$ws.Range("E3").Value = "  -19.84%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "439.05"
$ws.Range("E5").Value = "  -17.01%  "
$ws.Range("D6").Value = "122.00"
$ws.Range("E6").Value = "  -15.27%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "0.475"
$ws.Range("E8").Value = "  -14.49%  "
$ws.Range("D9").Value = "2.339.66"
$ws.Range("E9").Value = "  -19.55%  "
$ws.Range("D10").Value = "0.0920"
$ws.Range("E10").Value = "  -14.88%  "
$ws.Range("D11").Value = "5.30"
$ws.Range("E11").Value = "  -12.43%  "
$ws.Range("D12").Value = "0.310"
$ws.Range("E12").Value = "  -14.19%  "
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("D14").Value = "2.726.72"
$ws.Range("E14").Value = "  -20.18%  "
$ws.Range("D15").Value = "53.470.01"
$ws.Range("E15").Value = "  -11.86%  "
$ws.Range("D16").Value = "19.23"
$ws.Range("E16").Value = "  -15.89%  "
$ws.Range("D17").Value = "0.0000121"
$ws.Range("E17").Value = "  -14.83%  "
$ws.Range("D18").Value = "2.349.10"
$ws.Range("E18").Value = "  -19.31%  "
$ws.Range("D19").Value = "4.00"
$ws.Range("E19").Value = "  -20.68%  "
$ws.Range("D20").Value = "302.81"
$ws.Range("E20").Value = "  -16.54%  "
$ws.Range("D21").Value = "9.23"
$ws.Range("E21").Value = "  -21.46%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "5.58"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "5.41"
$ws.Range("E24").Value = "  -18.50%  "
$ws.Range("D25").Value = "55.63"
$ws.Range("E25").Value = "  -14.08%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "0.154"
$ws.Range("E27").Value = "  -14.09%  "
$ws.Range("D28").Value = "0.370"
$ws.Range("E28").Value = "  -18.80%  "
$ws.Range("D29").Value = "7.01"
$ws.Range("E29").Value = "  -10.86%  "
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").Value = "0.0₃0707"
$ws.Range("E31").Value = "  -17.35%  "
$ws.Range("D32").Value = "145.32"
$ws.Range("E32").Value = "  -4.66%  "
$ws.Range("D33").Value = "17.27"
$ws.Range("E33").Value = "  -12.67%  "
$ws.Range("E34").Value = "  -19.73%  "
$ws.Range("D35").Value = "4.74"
$ws.Range("E35").Value = "  -15.06%  "
$ws.Range("D36").Value = "3.57"
$ws.Range("E36").Value = "  -18.29%  "
$ws.Range("D37").Value = "0.829"
$ws.Range("E37").Value = "  -17.49%  "
$ws.Range("D38").Value = "1.01"
$ws.Range("E38").Value = "  -16.56%  "
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "33.10"
$ws.Range("E40").Value = "  -12.48%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "3.18"
$ws.Range("E42").Value = "  -14.69%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.22"
$ws.Range("E43").Value = "  -18.16%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "0.0500"
$ws.Range("E44").Value = "  -14.34%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.916.40"
$ws.Range("E45").Value = "  -16.61%  "
$ws.Range("D46").Value = "0.526"
$ws.Range("E46").Value = "  -19.14%  "
$ws.Range("D47").Value = "0.0210"
$ws.Range("E47").Value = "  -11.64%  "
$ws.Range("D48").Value = "0.0833"
$ws.Range("E48").Value = "  -9.88%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "3.97"
$ws.Range("E49").Value = "  -20.81%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "15.65"
$ws.Range("E50").Value = "  -23.66%  "
$ws.Range("E51").Value = "  -4.23%  "
